# "Discussed questions with Lukas"
# - The "Fragen an Lukas" sheet is replaced by a new "TODO" sheet with a
#   freshly curated list of questions/notes (some old rows dropped, several
#   new rows/answers added).
# - "Fragen an Ivan" gets one additional question appended at the bottom.

$wb = $excel.ActiveWorkbook

$wsIvan = $wb.Worksheets.Item("Fragen an Ivan")

# Recreate the second sheet from scratch (positioned right after "Fragen an
# Ivan") so it gets a fresh sheetId, then drop the old "Fragen an Lukas"
# sheet and rename the new one to "TODO".
$wsTodo = $wb.Worksheets.Add([System.Type]::Missing, $wsIvan)
$wsTodo.Name = "TODO_NEW_TMP"
$wb.Worksheets.Item("Fragen an Lukas").Delete() | Out-Null
$wsTodo.Name = "TODO"

# --- Populate "TODO" with the discussed / updated questions ---
$wsTodo.Range("A1").Value = "Findings "
$wsTodo.Range("B1").Value = "Finds"

$wsTodo.Range("A2").Value = "Schnitt/Quadrant übersetzig section"

$wsTodo.Range("A3").Value = "Raum"
$wsTodo.Range("B3").Value = "Room"

$wsTodo.Range("A4").Value = "positionDarueber/Darunter links uf Position"

$wsTodo.Range("A5").Value = "Sollten alle Vasen als Gefaesse bezeichnet werden? "
$wsTodo.Range("B5").Value = "Ja"

$wsTodo.Range("A6").Value = "PublikationNr lösche?"
$wsTodo.Range("B6").Value = "Delete"
$wsTodo.Range("C6").Value = 1690

$wsTodo.Range("A7").Value = "Bilder für Klassen?"
$wsTodo.Range("B7").Value = "Nein"

$wsTodo.Range("A8").Value = "Bild-Kartonage Nr?"

$wsTodo.Range("A9").Value = "Plan-Kampagne?"

$wsTodo.Range("A10").Value = "Tagebuch-Abhub?"

$wsTodo.Range("A11").Value = "Hat Zeichnung ein Bild?"

$wsTodo.Range("A12").Value = "Plan auch StillImageRepresentation"

$wsTodo.Range("A13").Value = "Persoenlichkeit-MaskenNr?"
$wsTodo.Range("B13").Value = 2705

$wsTodo.Range("A14").Value = "Lage alles 0-1 und denne mehreri Lage Objekt, Lage nur für d Abkürzig vom ganze, wie verlinke? -> Ivan "

$wsTodo.Range("A15").Value = "Jpeg 2000 conversion 1000 bilder"

$wsTodo.Columns.Item(1).AutoFit() | Out-Null

# --- "Fragen an Ivan": append the new question discussed with Lukas ---
$wsIvan.Range("A17").Value = "Lage Situation besprechen "

# --- Restore per-sheet selections; leave "TODO" as the active tab ---
$wsIvan.Range("A18").Select() | Out-Null
$wsTodo.Range("A40").Select() | Out-Null
